# CHP maps - "available fuel" feature
# - m_fuel becomes the active/selected sheet (was m_steam before)
# - W_el and Q_th: the B4:B10 "available fuel" column is no longer computed
#   from a formula (C-700 / C-2100); it is hard-set to 0 so the map can be
#   overridden/driven by an "available fuel" input later.
# - m_fuel and TIT: B4:B10 values reset to 0 for the same reason.
# - selections on every sheet are refreshed to reflect where the author was
#   last working.

$wb = $excel.ActiveWorkbook

$ws_m_steam = $wb.Worksheets.Item("m_steam")
$ws_w_el    = $wb.Worksheets.Item("W_el")
$ws_q_th    = $wb.Worksheets.Item("Q_th")
$ws_m_fuel  = $wb.Worksheets.Item("m_fuel")
$ws_tit     = $wb.Worksheets.Item("TIT")
$ws_tstack  = $wb.Worksheets.Item("Tstack")

# ---------------------------------------------------------------------
# W_el: replace the shared formula "=C{row}-700" in B4:B10 with a plain 0
# ---------------------------------------------------------------------
foreach ($r in 4..10) {
    $ws_w_el.Cells.Item($r, 2).Value = 0
}

# ---------------------------------------------------------------------
# Q_th: replace the shared formula "=C{row}-2100" in B4:B10 with a plain 0
# ---------------------------------------------------------------------
foreach ($r in 4..10) {
    $ws_q_th.Cells.Item($r, 2).Value = 0
}

# ---------------------------------------------------------------------
# m_fuel: B4:B10 hard values reset to 0
# ---------------------------------------------------------------------
foreach ($r in 4..10) {
    $ws_m_fuel.Cells.Item($r, 2).Value = 0
}

# ---------------------------------------------------------------------
# TIT: B4:B10 hard values reset to 0
# ---------------------------------------------------------------------
foreach ($r in 4..10) {
    $ws_tit.Cells.Item($r, 2).Value = 0
}

# ---------------------------------------------------------------------
# Selections (also drives which sheet is active/tabSelected: whichever
# sheet we select on last becomes the active tab)
# ---------------------------------------------------------------------
$ws_m_steam.Range("C14").Select()
$ws_w_el.Range("C14").Select()
$ws_q_th.Range("B3:B10").Select()
$ws_tit.Range("B3:B10").Select()
$ws_tstack.Range("J18").Select()

# m_fuel is the sheet that ends up active/selected in the saved workbook
$ws_m_fuel.Range("I3").Select()
